# Edit script for lab1.xlsx "end of first year" update
# Updates the random-sample data block (A2:J21), re-selects the active cell,
# and refreshes the historical P* log column (W2:W11) with this run's results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2:J21 sample data (columns X1..X10) ---
$ws.Cells.Item(2, 1).Value = 2.5252641376995149
$ws.Cells.Item(2, 2).Value = 6.0857267983031704
$ws.Cells.Item(2, 3).Value = 6.7977480391857661
$ws.Cells.Item(2, 4).Value = 5.839415570543534
$ws.Cells.Item(2, 5).Value = 7.1509961241492963
$ws.Cells.Item(2, 6).Value = 4.9145899838251896
$ws.Cells.Item(2, 7).Value = 6.7426755577257609
$ws.Cells.Item(2, 8).Value = 6.50064180425428
$ws.Cells.Item(2, 9).Value = 5.1054722739341409
$ws.Cells.Item(2, 10).Value = 6.960292062135685
$ws.Cells.Item(3, 1).Value = 2.3507788323618275
$ws.Cells.Item(3, 2).Value = 5.1441477706228831
$ws.Cells.Item(3, 3).Value = 4.7377877132480855
$ws.Cells.Item(3, 4).Value = 1.9488744773705253
$ws.Cells.Item(3, 5).Value = 6.4011905270546583
$ws.Cells.Item(3, 6).Value = 2.959606006042665
$ws.Cells.Item(3, 7).Value = 3.6226145207068088
$ws.Cells.Item(3, 8).Value = 2.4867668691061127
$ws.Cells.Item(3, 9).Value = 3.0469377727591782
$ws.Cells.Item(3, 10).Value = 2.0125019074068424
$ws.Cells.Item(4, 1).Value = 6.4257860042115542
$ws.Cells.Item(4, 2).Value = 5.7918286690878018
$ws.Cells.Item(4, 3).Value = 4.9488097781304363
$ws.Cells.Item(4, 4).Value = 6.8587020477919856
$ws.Cells.Item(4, 5).Value = 3.2699011200292976
$ws.Cells.Item(4, 6).Value = 3.8462907803582871
$ws.Cells.Item(4, 7).Value = 6.0267332987456887
$ws.Cells.Item(4, 8).Value = 6.4140229499191257
$ws.Cells.Item(4, 9).Value = 7.3335016937772757
$ws.Cells.Item(4, 10).Value = 2.7312958159123508
$ws.Cells.Item(5, 1).Value = 5.9484911648915064
$ws.Cells.Item(5, 2).Value = 3.4247813348796043
$ws.Cells.Item(5, 3).Value = 2.2142561113315224
$ws.Cells.Item(5, 4).Value = 3.1456761375774409
$ws.Cells.Item(5, 5).Value = 5.0614499343852044
$ws.Cells.Item(5, 6).Value = 4.1164845728934596
$ws.Cells.Item(5, 7).Value = 2.910949736014893
$ws.Cells.Item(5, 8).Value = 2.1515198217719047
$ws.Cells.Item(5, 9).Value = 5.7265971861934259
$ws.Cells.Item(5, 10).Value = 6.1659294412060905
$ws.Cells.Item(6, 1).Value = 5.0343592638935517
$ws.Cells.Item(6, 2).Value = 3.1893420209356975
$ws.Cells.Item(6, 3).Value = 6.3610892056031982
$ws.Cells.Item(6, 4).Value = 2.7008188116092411
$ws.Cells.Item(6, 5).Value = 2.741454817346721
$ws.Cells.Item(6, 6).Value = 6.806837672048097
$ws.Cells.Item(6, 7).Value = 6.6270055238502152
$ws.Cells.Item(6, 8).Value = 5.4264610736411631
$ws.Cells.Item(6, 9).Value = 1.8456804101687674
$ws.Cells.Item(6, 10).Value = 4.1786861781670579
$ws.Cells.Item(7, 1).Value = 4.2344715720084229
$ws.Cells.Item(7, 2).Value = 5.3029490035706655
$ws.Cells.Item(7, 3).Value = 4.1927661976989041
$ws.Cells.Item(7, 4).Value = 4.2678002258369698
$ws.Cells.Item(7, 5).Value = 1.7952418591875974
$ws.Cells.Item(7, 6).Value = 2.3566603595080418
$ws.Cells.Item(7, 7).Value = 6.5587441633350627
$ws.Cells.Item(7, 8).Value = 5.7424594866786709
$ws.Cells.Item(7, 9).Value = 4.1974001281777396
$ws.Cells.Item(7, 10).Value = 5.5266252632221438
$ws.Cells.Item(8, 1).Value = 1.8387295144505142
$ws.Cells.Item(8, 2).Value = 2.8970479445783868
$ws.Cells.Item(8, 3).Value = 1.5724567400128178
$ws.Cells.Item(8, 4).Value = 5.6497808771019624
$ws.Cells.Item(8, 5).Value = 2.4236741233558154
$ws.Cells.Item(8, 6).Value = 6.7777864925077056
$ws.Cells.Item(8, 7).Value = 5.9322724082155824
$ws.Cells.Item(8, 8).Value = 3.5263713492233038
$ws.Cells.Item(8, 9).Value = 3.4511590929898985
$ws.Cells.Item(8, 10).Value = 7.2335157322916341
$ws.Cells.Item(9, 1).Value = 2.2196029541917173
$ws.Cells.Item(9, 2).Value = 2.7295135349589525
$ws.Cells.Item(9, 3).Value = 5.1090368358409375
$ws.Cells.Item(9, 4).Value = 3.0991586046937467
$ws.Cells.Item(9, 5).Value = 4.4966451002533034
$ws.Cells.Item(9, 6).Value = 2.3801864680928984
$ws.Cells.Item(9, 7).Value = 5.8408413953062528
$ws.Cells.Item(9, 8).Value = 2.8679967650379954
$ws.Cells.Item(9, 9).Value = 2.0911004974517047
$ws.Cells.Item(9, 10).Value = 4.669348124637593
$ws.Cells.Item(10, 1).Value = 3.8664305551316875
$ws.Cells.Item(10, 2).Value = 6.2607467879268768
$ws.Cells.Item(10, 3).Value = 3.943425092318491
$ws.Cells.Item(10, 4).Value = 2.2727149266029847
$ws.Cells.Item(10, 5).Value = 6.4698083437604907
$ws.Cells.Item(10, 6).Value = 3.3693523972289192
$ws.Cells.Item(10, 7).Value = 6.8595931882686845
$ws.Cells.Item(10, 8).Value = 7.1581252479628894
$ws.Cells.Item(10, 9).Value = 2.7081261635181737
$ws.Cells.Item(10, 10).Value = 1.7841917172765283
$ws.Cells.Item(11, 1).Value = 2.4680529190954315
$ws.Cells.Item(11, 2).Value = 4.1676360362559892
$ws.Cells.Item(11, 3).Value = 4.1122070986053041
$ws.Cells.Item(11, 4).Value = 1.5785164952543718
$ws.Cells.Item(11, 5).Value = 6.2689452803125096
$ws.Cells.Item(11, 6).Value = 1.8396206549272134
$ws.Cells.Item(11, 7).Value = 3.9133045442060608
$ws.Cells.Item(11, 8).Value = 4.4813174840540784
$ws.Cells.Item(11, 9).Value = 1.6494512771996217
$ws.Cells.Item(11, 10).Value = 5.4505218665120392
$ws.Cells.Item(12, 1).Value = 6.0866179387798702
$ws.Cells.Item(12, 2).Value = 6.314928128910184
$ws.Cells.Item(12, 3).Value = 4.1795773186437577
$ws.Cells.Item(12, 4).Value = 2.979745780816065
$ws.Cells.Item(12, 5).Value = 3.9450291451765498
$ws.Cells.Item(12, 6).Value = 6.4628574480422367
$ws.Cells.Item(12, 7).Value = 4.7379659413434245
$ws.Cells.Item(12, 8).Value = 5.0618063905758843
$ws.Cells.Item(12, 9).Value = 7.0166121402630699
$ws.Cells.Item(12, 10).Value = 6.8105804620502335
$ws.Cells.Item(13, 1).Value = 2.6218637653736989
$ws.Cells.Item(13, 2).Value = 1.8643943601794488
$ws.Cells.Item(13, 3).Value = 1.7756367687002168
$ws.Cells.Item(13, 4).Value = 4.7192519913327438
$ws.Cells.Item(13, 5).Value = 3.9945765556810207
$ws.Cells.Item(13, 6).Value = 6.4398660237433996
$ws.Cells.Item(13, 7).Value = 3.3631144138920255
$ws.Cells.Item(13, 8).Value = 5.0372109134189884
$ws.Cells.Item(13, 9).Value = 4.1000875881221965
$ws.Cells.Item(13, 10).Value = 6.6074004333628347
$ws.Cells.Item(14, 1).Value = 2.1050022888882109
$ws.Cells.Item(14, 2).Value = 3.1510229804376353
$ws.Cells.Item(14, 3).Value = 3.8003079317606128
$ws.Cells.Item(14, 4).Value = 6.5929639576403085
$ws.Cells.Item(14, 5).Value = 2.8458964812158571
$ws.Cells.Item(14, 6).Value = 4.9595034638508251
$ws.Cells.Item(14, 7).Value = 2.9464171269875177
$ws.Cells.Item(14, 8).Value = 3.847360148930326
$ws.Cells.Item(14, 9).Value = 1.7856175420392468
$ws.Cells.Item(14, 10).Value = 6.5696160771507914
$ws.Cells.Item(15, 1).Value = 3.6153071687978757
$ws.Cells.Item(15, 2).Value = 6.3430881679738755
$ws.Cells.Item(15, 3).Value = 5.3950929288613549
$ws.Cells.Item(15, 4).Value = 5.4421451460310681
$ws.Cells.Item(15, 5).Value = 6.6059746086001159
$ws.Cells.Item(15, 6).Value = 5.8572383800775167
$ws.Cells.Item(15, 7).Value = 6.1803659169286167
$ws.Cells.Item(15, 8).Value = 6.6662157048249755
$ws.Cells.Item(15, 9).Value = 5.0272301400799586
$ws.Cells.Item(15, 10).Value = 6.6981185338908045
$ws.Cells.Item(16, 1).Value = 2.3518482009338664
$ws.Cells.Item(16, 2).Value = 2.6439640491958372
$ws.Cells.Item(16, 3).Value = 5.089253517258217
$ws.Cells.Item(16, 4).Value = 2.0208786278878139
$ws.Cells.Item(16, 5).Value = 5.3672893459883415
$ws.Cells.Item(16, 6).Value = 6.6733448286385695
$ws.Cells.Item(16, 7).Value = 2.0110760826441236
$ws.Cells.Item(16, 8).Value = 3.0344618060853907
$ws.Cells.Item(16, 9).Value = 5.6018375194555494
$ws.Cells.Item(16, 10).Value = 1.9037827692495499
$ws.Cells.Item(17, 1).Value = 5.3977663502914517
$ws.Cells.Item(17, 2).Value = 1.5610501419110692
$ws.Cells.Item(17, 3).Value = 3.6388332773827328
$ws.Cells.Item(17, 4).Value = 2.9551503036591695
$ws.Cells.Item(17, 5).Value = 6.9581533249916072
$ws.Cells.Item(17, 6).Value = 6.291223792229987
$ws.Cells.Item(17, 7).Value = 7.3545326090273742
$ws.Cells.Item(17, 8).Value = 4.0760267952513196
$ws.Cells.Item(17, 9).Value = 2.4368630024109623
$ws.Cells.Item(17, 10).Value = 3.9976064333017973
$ws.Cells.Item(18, 1).Value = 3.0057670827356793
$ws.Cells.Item(18, 2).Value = 5.2874431592761013
$ws.Cells.Item(18, 3).Value = 7.3541761528366951
$ws.Cells.Item(18, 4).Value = 3.5705719168675802
$ws.Cells.Item(18, 5).Value = 4.8967671742912078
$ws.Cells.Item(18, 6).Value = 7.1948402356028929
$ws.Cells.Item(18, 7).Value = 3.2208883938108461
$ws.Cells.Item(18, 8).Value = 2.4548640400402846
$ws.Cells.Item(18, 9).Value = 4.9554042176580095
$ws.Cells.Item(18, 10).Value = 2.3871373638111515
$ws.Cells.Item(19, 1).Value = 4.9498791467024752
$ws.Cells.Item(19, 2).Value = 4.153377788628803
$ws.Cells.Item(19, 3).Value = 2.7127600939970091
$ws.Cells.Item(19, 4).Value = 5.9508081301309241
$ws.Cells.Item(19, 5).Value = 5.0864018677327802
$ws.Cells.Item(19, 6).Value = 2.2996273689992979
$ws.Cells.Item(19, 7).Value = 7.0977059236426889
$ws.Cells.Item(19, 8).Value = 7.0144734031189913
$ws.Cells.Item(19, 9).Value = 3.0988021485030672
$ws.Cells.Item(19, 10).Value = 1.5596243171483506
$ws.Cells.Item(20, 1).Value = 3.6571907712027345
$ws.Cells.Item(20, 2).Value = 6.4952949613940856
$ws.Cells.Item(20, 3).Value = 3.8361317789239173
$ws.Cells.Item(20, 4).Value = 2.5489684743797114
$ws.Cells.Item(20, 5).Value = 6.3812289803765987
$ws.Cells.Item(20, 6).Value = 5.4412540055543692
$ws.Cells.Item(20, 7).Value = 6.5074144718771931
$ws.Cells.Item(20, 8).Value = 3.5030234687337867
$ws.Cells.Item(20, 9).Value = 6.6950886562700278
$ws.Cells.Item(20, 10).Value = 4.0164986114078189
$ws.Cells.Item(21, 1).Value = 2.6143781853694268
$ws.Cells.Item(21, 2).Value = 4.2882964568010502
$ws.Cells.Item(21, 3).Value = 3.0408780175176244
$ws.Cells.Item(21, 4).Value = 4.0125775933103425
$ws.Cells.Item(21, 5).Value = 4.1927661976989041
$ws.Cells.Item(21, 6).Value = 7.0923590807824946
$ws.Cells.Item(21, 7).Value = 5.3745966978972746
$ws.Cells.Item(21, 8).Value = 4.8443681142612993
$ws.Cells.Item(21, 9).Value = 3.9894079409161654
$ws.Cells.Item(21, 10).Value = 5.2184688863795889

# --- W2:W11 historical P* results log ---
$ws.Cells.Item(2, 23).Value = 0.615
$ws.Cells.Item(3, 23).Value = 0.54
$ws.Cells.Item(4, 23).Value = 0.62
$ws.Cells.Item(5, 23).Value = 0.65
$ws.Cells.Item(6, 23).Value = 0.575
$ws.Cells.Item(7, 23).Value = 0.585
$ws.Cells.Item(8, 23).Value = 0.59
$ws.Cells.Item(9, 23).Value = 0.595
$ws.Cells.Item(10, 23).Value = 0.6
$ws.Cells.Item(11, 23).Value = 0.565

# --- Active selection moves to AF34 after the edit ---
$ws.Range("AF34").Select()
